$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "Implement Stack using Array / Linked List" row was split into three
# separate, completed problems: "Implement Stack using Array",
# "Implement Stack using Linked List" and "Implement Stack using Queues"
# (previously further down, blank). All the following Stack rows shift up
# by two positions, and a couple of bookkeeping rows below get tidied up.
# ---------------------------------------------------------------------------

function Set-DateCell($addr, $value) {
    # Give a brand-new date cell the same number format as an existing,
    # already-formatted date cell (column F date style) by copying formats.
    $ws.Range($addr).Value = $value
    $ws.Range("F110").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

# Row 111: Implement Stack using Array (new, completed)
$ws.Range("C111").Value = "Implement Stack using Array"
$ws.Range("D111").Value = "Easy"
$ws.Range("E111").Value = "Done"
Set-DateCell "F111" 45665
$ws.Range("G111").Value = "O(1)"
$ws.Range("H111").Value = "O(n)"

# Row 112: Implement Stack using Linked List (new, completed)
$ws.Range("C112").Value = "Implement Stack using Linked List"
$ws.Range("D112").Value = "Easy"
$ws.Range("E112").Value = "Done"
Set-DateCell "F112" 45665
$ws.Range("G112").Value = "O(1)"
$ws.Range("H112").Value = "O(1)"

# Row 113: Implement Stack using Queues (completed)
$ws.Range("C113").Value = "Implement Stack using Queues"
$ws.Range("F113").Value = 45665
$ws.Range("G113").Value = "Push O(n) Pop O(1)"
$ws.Range("H113").Value = "O(1)"
$ws.Range("I113").Clear() | Out-Null

# Row 114: Next Greater Element I
$ws.Range("C114").Value = "Next Greater Element I"
$ws.Range("D114").Value = "Easy"
$ws.Range("G114").Value = "O(m + n)"
$ws.Range("I114").Value = "Monotonic Stack"

# Row 115: Next Greater Element II (Circular)
$ws.Range("C115").Value = "Next Greater Element II (Circular)"
$ws.Range("I115").Value = "Monotonic Stack"

# Row 116: Daily Temperatures
$ws.Range("C116").Value = "Daily Temperatures"
$ws.Range("I116").Value = "Monotonic Stack"

# Row 117: Evaluate Reverse Polish Notation (now completed)
$ws.Range("C117").Value = "Evaluate Reverse Polish Notation"
$ws.Range("D117").Value = "Medium"
$ws.Range("E117").Value = "Done"
Set-DateCell "F117" 45664
$ws.Range("G117").Value = "O(n)"
$ws.Range("H117").Value = "O(n)"

# Row 118: Largest Rectangle in Histogram
$ws.Range("C118").Value = "Largest Rectangle in Histogram"

# Row 119: Trapping Rain Water (Stack approach)
$ws.Range("C119").Value = "Trapping Rain Water (Stack approach)"

# Row 120: Remove K Digits
$ws.Range("C120").Value = "Remove K Digits"

# Row 121: Decode String (now To Do / not yet attempted)
$ws.Range("C121").Value = "Decode String"
$ws.Range("D121:H121").Clear() | Out-Null

# Row 122: Asteroid Collision
$ws.Range("C122").Value = "Asteroid Collision"
$ws.Range("G122").Value = "O(n)"
$ws.Range("I122").Clear() | Out-Null

# Row 123: Online Stock Span
$ws.Range("C123").Value = "Online Stock Span"
$ws.Range("F123").Value = 45665
$ws.Range("G123").Value = "O(1)"
$ws.Range("I123").Value = "Monotonic Stack"

# Row 124: Simplify Path (now filled in)
$ws.Range("B124").Value = "Stack"
$ws.Range("C124").Value = "Simplify Path"
$ws.Range("D124").Value = "Medium"
$ws.Range("E124").Value = "Done"
Set-DateCell "F124" 45662
$ws.Range("G124").Value = "O(n)"
$ws.Range("H124").Value = "O(n)"

# Row 126 loses its contents (the "Valid Anagram" placeholder moves to 127)
$ws.Range("B126:E126").Clear() | Out-Null

# Row 127: String / Valid Anagram (replaces the old Stack / Valid Parentheses dup)
$ws.Range("B127").Value = "String"
$ws.Range("C127").Value = "Valid Anagram"

# Column G is slightly wider to fit the new text.
$ws.Columns("G:G").ColumnWidth = 15.5

# Match the author's final on-screen selection/scroll position.
$ws.Range("G114").Select() | Out-Null

Write-Output "done"
